$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 21198.785
$ws.Range("J19").Value = 28098.8
$ws.Range("L19").Value = 28098.8
$ws.Range("N19").Value = -28448.8

$ws.Range("H40").Value = 2700
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 2700
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 2700
$ws.Range("N40").Value = -3050
$ws.Range("M40").ClearContents()

$ws.Range("H92").Value = 33823.535
$ws.Range("I92").Value = 480.76
$ws.Range("K92").Value = 480.76
$ws.Range("M92").Value = 767.24

$ws.Range("H99").Value = 277.375
$ws.Range("I99").Value = 263.16666
$ws.Range("J99").Value = 320
$ws.Range("K99").Value = 789.4999799999999
$ws.Range("L99").Value = 960
$ws.Range("M99").Value = 708.5000200000001
$ws.Range("N99").Value = -3956

$ws.Range("H100").Value = 3793.1667
$ws.Range("I100").Value = 4028.5
$ws.Range("K100").Value = 4028.5
$ws.Range("M100").Value = -3487.5

$ws.Range("H101").Value = 336.5
$ws.Range("I101").Value = 348.08334
$ws.Range("K101").Value = 1044.25002
$ws.Range("M101").Value = 577.7499800000001

$ws.Range("H111").Value = 83748.25
$ws.Range("I111").Value = 4900
$ws.Range("J111").Value = 110031
$ws.Range("K111").Value = 14700
$ws.Range("L111").Value = 330093
$ws.Range("M111").Value = -11633
$ws.Range("N111").Value = -336227

$ws.Range("H137").Value = 3157.125
$ws.Range("I137").Value = 2351.6956
$ws.Range("J137").Value = 5215.4443
$ws.Range("K137").Value = 7055.0868
$ws.Range("L137").Value = 15646.3329
$ws.Range("M137").Value = -4505.0868
$ws.Range("N137").Value = -20746.3329

$ws.Range("H141").Value = 5368.4287
$ws.Range("I141").Value = 6233.4546
$ws.Range("K141").Value = 18700.3638
$ws.Range("M141").Value = -13520.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6387.2
$ws.Range("I32").Value = 4872
$ws.Range("J32").Value = 18761.334
$ws.Range("K32").Value = 4872
$ws.Range("L32").Value = 18761.334
$ws.Range("M32").Value = -4585
$ws.Range("N32").Value = -19335.334

$ws.Range("H35").Value = 1212.4286
$ws.Range("I35").Value = 1212.4286
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1212.4286
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -806.4286
$ws.Range("N35").ClearContents()

$ws.Range("H45").Value = 3005.125
$ws.Range("I45").Value = 2924.5
$ws.Range("J45").Value = 3247
$ws.Range("K45").Value = 2924.5
$ws.Range("L45").Value = 3247
$ws.Range("M45").Value = -2547.5
$ws.Range("N45").Value = -4001

$ws.Range("H74").Value = 1643.7736
$ws.Range("I74").Value = 1534.2084
$ws.Range("K74").Value = 1534.2084
$ws.Range("M74").Value = -660.2084

$ws.Range("H77").Value = 1643.7736
$ws.Range("I77").Value = 1534.2084
$ws.Range("K77").Value = 7671.041999999999
$ws.Range("M77").Value = -3303.041999999999

$ws.Range("H97").Value = 1532.8462
$ws.Range("I97").Value = 1621.8
$ws.Range("K97").Value = 1621.8
$ws.Range("M97").Value = -1125.8

$ws.Range("H102").Value = 5982474
$ws.Range("I102").Value = 6765684.5
$ws.Range("K102").Value = 6765684.5
$ws.Range("M102").Value = -6764062.5

$ws.Range("H122").Value = 2185.3403
$ws.Range("I122").Value = 1574.7941
$ws.Range("J122").Value = 3782.1538
$ws.Range("K122").Value = 4724.3823
$ws.Range("L122").Value = 11346.4614
$ws.Range("M122").Value = -2274.3823
$ws.Range("N122").Value = -16246.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 41667092
$ws.Range("J80").Value = 71428930
$ws.Range("L80").Value = 71428930
$ws.Range("N80").Value = -71430926

$ws.Range("H83").Value = 41667092
$ws.Range("J83").Value = 71428930
$ws.Range("L83").Value = 357144650
$ws.Range("N83").Value = -357154634

$ws.Range("H105").Value = 2082.0312
$ws.Range("I105").Value = 1830.9546
$ws.Range("K105").Value = 1830.9546
$ws.Range("M105").Value = -83.95460000000003

$ws.Range("H134").Value = 2546.2954
$ws.Range("I134").Value = 2312.4473
$ws.Range("K134").Value = 6937.341899999999
$ws.Range("M134").Value = -4402.341899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.714287
$ws.Range("I7").Value = 14.666667
$ws.Range("J7").Value = 90
$ws.Range("K7").Value = 14.666667
$ws.Range("L7").Value = 90
$ws.Range("M7").Value = 98.333333
$ws.Range("N7").Value = -316

$ws.Range("H22").Value = 633.7368
$ws.Range("I22").Value = 375.14285
$ws.Range("K22").Value = 375.14285
$ws.Range("M22").Value = -25.14285000000001

$ws.Range("H31").Value = 2426.303
$ws.Range("J31").Value = 2983.2307
$ws.Range("L31").Value = 2983.2307
$ws.Range("N31").Value = -3573.2307

$ws.Range("H34").Value = 2426.303
$ws.Range("J34").Value = 2983.2307
$ws.Range("L34").Value = 2983.2307
$ws.Range("N34").Value = -3387.2307

$ws.Range("H107").Value = 556.11536
$ws.Range("I107").Value = 512.8095
$ws.Range("K107").Value = 512.8095
$ws.Range("M107").Value = 1407.1905

$ws.Range("H122").Value = 3357
$ws.Range("I122").Value = 5197.6
$ws.Range("J122").Value = 1823.1666
$ws.Range("K122").Value = 15592.8
$ws.Range("L122").Value = 5469.4998
$ws.Range("M122").Value = -13142.8
$ws.Range("N122").Value = -10369.4998

$ws.Range("H134").Value = 2719.9333
$ws.Range("I134").Value = 2272.182
$ws.Range("K134").Value = 6816.545999999999
$ws.Range("M134").Value = -4281.545999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 302.64706
$ws.Range("I38").Value = 176
$ws.Range("J38").Value = 483.57144
$ws.Range("K38").Value = 528
$ws.Range("L38").Value = 1450.71432
$ws.Range("M38").Value = -181
$ws.Range("N38").Value = -2144.71432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8297.467000000001
$ws.Range("I70").Value = 6947.6
$ws.Range("J70").Value = 10997.2
$ws.Range("K70").Value = 6947.6
$ws.Range("L70").Value = 10997.2
$ws.Range("M70").Value = -6677.6
$ws.Range("N70").Value = -11537.2

$ws.Range("H73").Value = 8297.467000000001
$ws.Range("I73").Value = 6947.6
$ws.Range("J73").Value = 10997.2
$ws.Range("K73").Value = 6947.6
$ws.Range("L73").Value = 10997.2
$ws.Range("M73").Value = -6011.6
$ws.Range("N73").Value = -12869.2

$ws.Range("H132").Value = 4349.457
$ws.Range("I132").Value = 4477.9585
$ws.Range("K132").Value = 13433.8755
$ws.Range("M132").Value = -10903.8755

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22730670
$ws.Range("J7").Value = 4202.5
$ws.Range("L7").Value = 4202.5
$ws.Range("N7").Value = -4426.5

$ws.Range("H22").Value = 764.4
$ws.Range("I22").Value = 875
$ws.Range("J22").Value = 690.6667
$ws.Range("K22").Value = 875
$ws.Range("L22").Value = 690.6667
$ws.Range("M22").Value = -580
$ws.Range("N22").Value = -1280.6667

$ws.Range("H27").Value = 764.4
$ws.Range("I27").Value = 875
$ws.Range("J27").Value = 690.6667
$ws.Range("K27").Value = 875
$ws.Range("L27").Value = 690.6667
$ws.Range("M27").Value = -768
$ws.Range("N27").Value = -904.6667

$ws.Range("H55").Value = 481.22726
$ws.Range("I55").Value = 183.92308
$ws.Range("J55").Value = 910.6667
$ws.Range("K55").Value = 183.92308
$ws.Range("L55").Value = 910.6667
$ws.Range("M55").Value = -10.92308
$ws.Range("N55").Value = -1256.6667

$ws.Range("H56").Value = 9480.799999999999
$ws.Range("I56").Value = 8112.25
$ws.Range("K56").Value = 8112.25
$ws.Range("M56").Value = -7421.25

$ws.Range("H93").Value = 1990.5652
$ws.Range("I93").Value = 1894.238
$ws.Range("K93").Value = 1894.238
$ws.Range("M93").Value = -646.2380000000001

$ws.Range("H126").Value = 22730670
$ws.Range("J126").Value = 4202.5
$ws.Range("L126").Value = 12607.5
$ws.Range("N126").Value = -17547.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 179000.38
$ws.Range("I4").Value = 32001
$ws.Range("J4").Value = 325999.75
$ws.Range("K4").Value = 32001
$ws.Range("L4").Value = 325999.75
$ws.Range("M4").Value = -31888
$ws.Range("N4").Value = -326225.75

$ws.Range("H34").Value = 5008000
$ws.Range("I34").Value = 5008000
$ws.Range("K34").Value = 5008000
$ws.Range("M34").Value = -5007797

$ws.Range("H58").Value = 18521.25
$ws.Range("I58").Value = 14085
$ws.Range("K58").Value = 14085
$ws.Range("M58").Value = -13777

$ws.Range("H113").Value = 1644.9565
$ws.Range("I113").Value = 544.3077
$ws.Range("J113").Value = 3075.8
$ws.Range("K113").Value = 1632.9231
$ws.Range("L113").Value = 9227.400000000001
$ws.Range("M113").Value = 537.0769
$ws.Range("N113").Value = -13567.4

$ws.Range("H126").Value = 1661.8125
$ws.Range("J126").Value = 3125
$ws.Range("L126").Value = 9375
$ws.Range("N126").Value = -14315

$ws.Range("H130").Value = 40808
$ws.Range("J130").Value = 40808
$ws.Range("L130").Value = 40808
$ws.Range("N130").Value = -50848

$ws.Range("H131").Value = 49715
$ws.Range("J131").Value = 49715
$ws.Range("L131").Value = 49715
$ws.Range("N131").Value = -59795

$ws.Range("H132").Value = 5131.3687
$ws.Range("I132").Value = 4963.8125
$ws.Range("J132").Value = 6025
$ws.Range("K132").Value = 14891.4375
$ws.Range("L132").Value = 18075
$ws.Range("M132").Value = -12361.4375
$ws.Range("N132").Value = -23135
